{"js": "// Replace each old division expression with its corresponding new one.\n// The mapping is a 1:1, collision-free set of exact-text replacements\n// (every \"old\" value is unique, every \"new\" value is unique, and no\n// \"old\" value equals any \"new\" value), so a straightforward\n// search-and-replace per pair is safe regardless of order.\nconst replacements = [\n  [\"856\u00f78=\", \"504\u00f72=\"],\n  [\"673\u00f78=\", \"133\u00f78=\"],\n  [\"142\u00f74=\", \"570\u00f79=\"],\n  [\"821\u00f75=\", \"474\u00f73=\"],\n  [\"719\u00f73=\", \"433\u00f79=\"],\n  [\"655\u00f79=\", \"572\u00f78=\"],\n  [\"947\u00f73=\", \"459\u00f75=\"],\n  [\"680\u00f75=\", \"367\u00f75=\"],\n  [\"603\u00f76=\", \"455\u00f74=\"],\n  [\"628\u00f73=\", \"292\u00f79=\"],\n  [\"935\u00f76=\", \"435\u00f74=\"],\n  [\"161\u00f79=\", \"286\u00f79=\"],\n  [\"510\u00f79=\", \"996\u00f79=\"],\n  [\"765\u00f75=\", \"668\u00f72=\"],\n  [\"318\u00f76=\", \"677\u00f76=\"],\n  [\"948\u00f72=\", \"631\u00f72=\"],\n  [\"270\u00f76=\", \"909\u00f76=\"],\n  [\"194\u00f73=\", \"999\u00f75=\"],\n  [\"731\u00f75=\", \"255\u00f72=\"],\n  [\"625\u00f78=\", \"177\u00f78=\"],\n  [\"812\u00f72=\", \"218\u00f73=\"],\n  [\"389\u00f76=\", \"800\u00f76=\"],\n  [\"394\u00f75=\", \"222\u00f72=\"],\n  [\"729\u00f79=\", \"170\u00f75=\"],\n  [\"899\u00f74=\", \"165\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old division expression with its corresponding new one.\n# The mapping is a 1:1, collision-free set of exact-text replacements\n# (every \"old\" value is unique, every \"new\" value is unique, and no\n# \"old\" value equals any \"new\" value), so a straightforward\n# Find/Replace per pair is safe regardless of order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"856\u00f78=\", \"504\u00f72=\"),\n    @(\"673\u00f78=\", \"133\u00f78=\"),\n    @(\"142\u00f74=\", \"570\u00f79=\"),\n    @(\"821\u00f75=\", \"474\u00f73=\"),\n    @(\"719\u00f73=\", \"433\u00f79=\"),\n    @(\"655\u00f79=\", \"572\u00f78=\"),\n    @(\"947\u00f73=\", \"459\u00f75=\"),\n    @(\"680\u00f75=\", \"367\u00f75=\"),\n    @(\"603\u00f76=\", \"455\u00f74=\"),\n    @(\"628\u00f73=\", \"292\u00f79=\"),\n    @(\"935\u00f76=\", \"435\u00f74=\"),\n    @(\"161\u00f79=\", \"286\u00f79=\"),\n    @(\"510\u00f79=\", \"996\u00f79=\"),\n    @(\"765\u00f75=\", \"668\u00f72=\"),\n    @(\"318\u00f76=\", \"677\u00f76=\"),\n    @(\"948\u00f72=\", \"631\u00f72=\"),\n    @(\"270\u00f76=\", \"909\u00f76=\"),\n    @(\"194\u00f73=\", \"999\u00f75=\"),\n    @(\"731\u00f75=\", \"255\u00f72=\"),\n    @(\"625\u00f78=\", \"177\u00f78=\"),\n    @(\"812\u00f72=\", \"218\u00f73=\"),\n    @(\"389\u00f76=\", \"800\u00f76=\"),\n    @(\"394\u00f75=\", \"222\u00f72=\"),\n    @(\"729\u00f79=\", \"170\u00f75=\"),\n    @(\"899\u00f74=\", \"165\u00f76=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
